$wb = $excel.ActiveWorkbook

# ---- Sheet "Fig.7" (sheet1): update Block Number column (A) values,
#      add new column F "IoVBlockSecure model" with data ----
$ws1 = $wb.Worksheets.Item("Fig.7")

# Update column A (Block Number) values: 10..100 step 10 -> 20..200 step 20
$ws1.Range("A2").Value = 20
$ws1.Range("A3").Value = 40
$ws1.Range("A4").Value = 60
$ws1.Range("A5").Value = 80
$ws1.Range("A6").Value = 100
$ws1.Range("A7").Value = 120
$ws1.Range("A8").Value = 140
$ws1.Range("A9").Value = 160
$ws1.Range("A10").Value = 180
$ws1.Range("A11").Value = 200

# New column F header + values
$ws1.Range("F1").Value = "IoVBlockSecure model"
$ws1.Range("F2").Value = 224
$ws1.Range("F3").Value = 448
$ws1.Range("F4").Value = 672
$ws1.Range("F5").Value = 896
$ws1.Range("F6").Value = 1120
$ws1.Range("F7").Value = 1344
$ws1.Range("F8").Value = 1568
$ws1.Range("F9").Value = 1792
$ws1.Range("F10").Value = 2016
$ws1.Range("F11").Value = 2240

# Give the new header cell F1 the same style as the other header cells
$ws1.Range("E1").Copy() | Out-Null
$ws1.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update selection on Fig.7 and make it the active sheet/tab
$ws1.Activate() | Out-Null
$ws1.Range("M10").Select() | Out-Null

# ---- Sheet "Fig.8" (sheet2): change selection only ----
$ws2 = $wb.Worksheets.Item("Fig.8")
$ws2.Activate() | Out-Null
$ws2.Range("A1:D1").Select() | Out-Null

# Restore Fig.7 as the active/selected sheet (matches target workbook state)
$ws1.Activate() | Out-Null
$ws1.Range("M10").Select() | Out-Null
